$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dces")
$ws.Activate()

# Add the missing row of date-item data (row 14) to the "dces" sheet
$ws.Range("A14").Value = 9
$ws.Range("D14").Value = 222
$ws.Range("J14").Value = "adfasdfas"

# Leave the selection where the user clicked after entering the data
$ws.Range("G13").Select()
